$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.649.32"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "2.473.22"
$ws.Range("E3").Value = "  -0.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.55"
$ws.Range("E5").Value = "  +1.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.56"
$ws.Range("E6").Value = "  +0.99%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  +0.82%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.02"
$ws.Range("E10").Value = "  +1.82%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  +8.56%  "

# Row 12
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("D13").Value = "2.854.98"

# Row 14
$ws.Range("E14").Value = "  +0.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.76"
$ws.Range("E15").Value = "  -3.09%  "

# Row 16
$ws.Range("D16").Value = "2.466.46"
$ws.Range("E16").Value = "  +0.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  +2.58%  "

# Row 18
$ws.Range("D18").Value = "41.598.26"
$ws.Range("E18").Value = "  +0.19%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0949"
$ws.Range("E19").Value = "  +0.58%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  -0.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.24"
$ws.Range("E21").Value = "  -0.76%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.30"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.65"
$ws.Range("E23").Value = "  +1.57%  "

# Row 24
$ws.Range("E24").Value = "  +0.75%  "

# Row 25
$ws.Range("E25").Value = "  +1.61%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.73"
$ws.Range("E27").Value = "  -0.39%  "

# Row 28
$ws.Range("E28").Value = "  +2.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.83"
$ws.Range("E29").Value = "  +1.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.14"
$ws.Range("E30").Value = "  +1.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.14"
$ws.Range("E31").Value = "  +2.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("E32").Value = "  +1.34%  "

# Row 33
$ws.Range("E33").Value = "  -0.04%  "

# Row 34
$ws.Range("E34").Value = "  +0.52%  "

# Row 35
$ws.Range("E35").Value = "  +1.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.30"
$ws.Range("E36").Value = "  +0.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.93"
$ws.Range("E37").Value = "  +1.40%  "

# Row 38
$ws.Range("E38").Value = "  +2.36%  "

# Row 39
$ws.Range("E39").Value = "  +1.72%  "

# Row 40
$ws.Range("E40").Value = "  -0.44%  "

# Row 41
$ws.Range("E41").Value = "  -1.25%  "

# Row 42
$ws.Range("E42").Value = "  +3.51%  "

# Row 43
$ws.Range("D43").Value = "1.984.78"
$ws.Range("E43").Value = "  +1.37%  "

# Row 44
$ws.Range("E44").Value = "  +0.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.73"
$ws.Range("E45").Value = "  +0.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +1.86%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.29"
$ws.Range("E47").Value = "  +2.78%  "

# Row 48
$ws.Range("D48").Value = "2.713.51"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.18"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.00"
$ws.Range("E50").Value = "  +2.77%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.13"
$ws.Range("E51").Value = "  -0.20%  "
